$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "type"
$ws.Range("C2").Value = "title"
$ws.Range("C3").Value = "button"
$ws.Range("C4").Value = "description"

[void]$ws.Range("C5").Select()
